$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "d=6" data column belongs between the existing "d=5" (F) and
# "d=7" (old G) columns. Insert a blank column at G so everything from
# the old G (d=7) onward shifts one column to the right (G->H, H->I).
$ws.Columns("G:G").Insert()

# Header for the newly inserted column.
$ws.Range("G1").Value = "d=6"

# Data for the newly inserted column.
$ws.Range("G2").Value = 97.70571711987347
$ws.Range("G3").Value = 97.83009366944107
$ws.Range("G4").Value = 97.75969875795748
$ws.Range("G5").Value = 97.70414593527714
$ws.Range("G6").Value = 97.75012670431555
